# Edit script for DELAWARE_2019.xlsx
# 1. Rename header columns to snake_case English names
# 2. Capitalize connector words (de/del/el/y/la/las/los) in municipality/state names
# 3. Adjust two rounding values in column D (rows 391 and 446)
# 4. Remove trailing metadata/footer rows (449-453)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Capitalize connector words in municipality / state names ---
$ws.Range("B13").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Mazapa De Madero"
$ws.Range("B28").Value = "San Cristóbal De Las Casas"
$ws.Range("B39").Value = "Coyame Del Sotol"
$ws.Range("B43").Value = "Guadalupe Y Calvo"
$ws.Range("A51").Value = "Ciudad De México"
$ws.Range("A71").Value = "Estado De México"
$ws.Range("B71").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B72").Value = "Almoloya De Alquisiras"
$ws.Range("B73").Value = "Almoloya De Juárez"
$ws.Range("B83").Value = "Ecatepec De Morelos"
$ws.Range("B86").Value = "Ixtapan De La Sal"
$ws.Range("B91").Value = "Naucalpan De Juárez"
$ws.Range("B96").Value = "San Felipe Del Progreso"
$ws.Range("B103").Value = "Tenango Del Valle"
$ws.Range("B106").Value = "Tlalnepantla De Baz"
$ws.Range("B111").Value = "Valle De Bravo"
$ws.Range("B112").Value = "Villa De Allende"
$ws.Range("B122").Value = "Apaseo El Grande"
$ws.Range("B137").Value = "Valle De Santiago"
$ws.Range("B140").Value = "Acapulco De Juárez"
$ws.Range("B143").Value = "Atoyac De Álvarez"
$ws.Range("B144").Value = "Ayutla De Los Libres"
$ws.Range("B146").Value = "Chilapa De Álvarez"
$ws.Range("B147").Value = "Chilpancingo De Los Bravo"
$ws.Range("B148").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B153").Value = "Huitzuco De Los Figueroa"
$ws.Range("B154").Value = "Iguala De La Independencia"
$ws.Range("B155").Value = "Zihuatanejo De Azueta"
$ws.Range("B165").Value = "Taxco De Alarcón"
$ws.Range("B167").Value = "Técpan De Galeana"
$ws.Range("B169").Value = "Tepecoacuilco De Trujano"
$ws.Range("B172").Value = "Tlapa De Comonfort"
$ws.Range("B178").Value = "Huejutla De Reyes"
$ws.Range("B180").Value = "Jacala De Ledezma"
$ws.Range("B181").Value = "Mineral Del Monte"
$ws.Range("B182").Value = "Pachuca De Soto"
$ws.Range("B184").Value = "Santiago De Anaya"
$ws.Range("B186").Value = "Tulancingo De Bravo"
$ws.Range("B189").Value = "Encarnación De Díaz"
$ws.Range("B191").Value = "Lagos De Moreno"
$ws.Range("B193").Value = "San Martín De Bolaños"
$ws.Range("B195").Value = "Tepatitlán De Morelos"
$ws.Range("B196").Value = "Zapotlán Del Rey"
$ws.Range("B197").Value = "Zapotlán El Grande"
$ws.Range("B225").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B237").Value = "Coatlán Del Río"
$ws.Range("B248").Value = "Tlaltizapán De Zapata"
$ws.Range("B267").Value = "Ixtlán De Juárez"
$ws.Range("B268").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B271").Value = "Oaxaca De Juárez"
$ws.Range("B272").Value = "Putla Villa De Guerrero"
$ws.Range("B292").Value = "Tataltepec De Valdés"
$ws.Range("B310").Value = "Huehuetlán El Chico"
$ws.Range("B312").Value = "Ixcamilpa De Guerrero"
$ws.Range("B314").Value = "Izúcar De Matamoros"
$ws.Range("B317").Value = "Los Reyes De Juárez"
$ws.Range("B319").Value = "Palmar De Bravo"
$ws.Range("B334").Value = "Tetela De Ocampo"
$ws.Range("B337").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B349").Value = "Amealco De Bonfil"
$ws.Range("B350").Value = "Cadereyta De Montes"
$ws.Range("B352").Value = "Pinal De Amoles"
$ws.Range("B359").Value = "Ciudad Del Maíz"
$ws.Range("B364").Value = "San Ciro De Acosta"
$ws.Range("B367").Value = "Villa De Reyes"
$ws.Range("B404").Value = "Cosamaloapan De Carpio"
$ws.Range("B408").Value = "Ignacio De La Llave"
$ws.Range("B409").Value = "Ixhuatlán De Madero"
$ws.Range("B413").Value = "Martínez De La Torre"
$ws.Range("B414").Value = "Medellín De Bravo"
$ws.Range("B421").Value = "Poza Rica De Hidalgo"
$ws.Range("B424").Value = "Soledad De Doblado"
$ws.Range("B443").Value = "Villa De Cos"

# --- Fix rounding of two percentage values ---
$ws.Range("D391").Value = 0.009566326530612243
$ws.Range("D446").Value = 0.009566326530612243

# --- Remove trailing metadata/footer rows (449-453) ---
$ws.Range("A449:A453").EntireRow.Delete()
